# Applies the DaCapo Shenandoah GC "jme / heap-4G" benchmark table updates.
#
# The document is a single table, one column, 46 rows. Most rows hold a
# single numeric/text value; the last three rows originally held several
# tab-separated values crammed into one cell (an old "raw line paste").
# This edit:
#   1. Updates the first block of summary values (rows 1-12).
#   2. Collapses the garbled multi-value rows 44-46 down to the single
#      correct summary value each (matching rows 1-3's new column header
#      semantics), removing the stray tab-separated data.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cell = $table.Cell($rowIndex, 1)
    $cell.Range.Text = $newText
}

# --- Rows 1-4: top summary values ---
Set-CellText $t 1 "0M"
Set-CellText $t 2 "0M"
Set-CellText $t 3 "0M"
Set-CellText $t 4 "80"

# --- Rows 6-12: per-phase timing values ---
Set-CellText $t 6  "0.00069"
Set-CellText $t 7  "0.00026"
Set-CellText $t 8  "0.00010"
Set-CellText $t 9  "0.00042"
Set-CellText $t 10 "0.00052"
Set-CellText $t 11 "0.00054"
Set-CellText $t 12 "0.02089"

# --- Rows 44-46: collapse stray tab-separated raw data into single values ---
Set-CellText $t 44 "99.97"
Set-CellText $t 45 "0.02"
Set-CellText $t 46 "70"
